$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 164.7314145
$ws.Range("H2").Value = 329.462829
$ws.Range("I2").Value = 0.1559412251502966
$ws.Range("J2").Value = 0.1160318088037158
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 67.073376
$ws.Range("N2").Value = 134.146752
$ws.Range("O2").Value = 0.3886152607163267
$ws.Range("P2").Value = 0.3201976973989085
$ws.Range("Q2").Value = 11049.09210377035
$ws.Range("R2").Value = 44196.36841508141
$ws.Range("S2").Value = 0.06060113986820593
$ws.Range("T2").Value = 0.0371531180039802

$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 164.7314145
$ws.Range("H3").Value = 329.462829
$ws.Range("I3").Value = 0.1559412251502966
$ws.Range("J3").Value = 0.1160318088037158
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 4.219010666666667
$ws.Range("N3").Value = 12.657032
$ws.Range("O3").Value = 0.02444445215030819
$ws.Range("P3").Value = 0.03021133528677833
$ws.Range("Q3").Value = 695.003594910588
$ws.Range("R3").Value = 4170.021569463528
$ws.Range("S3").Value = 0.003811897816446863
$ws.Range("T3").Value = 0.003505475879700415

$ws.Range("E4").Value = 2
$ws.Range("G4").Value = 164.7314145
$ws.Range("H4").Value = 329.462829
$ws.Range("I4").Value = 0.1559412251502966
$ws.Range("J4").Value = 0.1160318088037158
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 12.54988666666667
$ws.Range("N4").Value = 37.64966
$ws.Range("O4").Value = 0.07271256897710082
$ws.Range("P4").Value = 0.0898667635266472
$ws.Range("Q4").Value = 2067.36058241469
$ws.Range("R4").Value = 12404.16349448814
$ws.Range("S4").Value = 0.01133888709011455
$ws.Range("T4").Value = 0.01042740312333267

$ws.Range("E5").Value = 2
$ws.Range("G5").Value = 164.7314145
$ws.Range("H5").Value = 329.462829
$ws.Range("I5").Value = 0.1559412251502966
$ws.Range("J5").Value = 0.1160318088037158
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 12.18648933333333
$ws.Range("N5").Value = 36.559468
$ws.Range("O5").Value = 0.07060708751994335
$ws.Range("P5").Value = 0.0872645613643264
$ws.Range("Q5").Value = 2007.497625669162
$ws.Range("R5").Value = 12044.98575401497
$ws.Range("S5").Value = 0.01101055573215418
$ws.Range("T5").Value = 0.01012546489956565

$ws.Range("E6").Value = 2
$ws.Range("G6").Value = 164.7314145
$ws.Range("H6").Value = 329.462829
$ws.Range("I6").Value = 0.1559412251502966
$ws.Range("J6").Value = 0.1160318088037158
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 44.80271133333334
$ws.Range("N6").Value = 134.408134
$ws.Range("O6").Value = 0.2595816460111038
$ws.Range("P6").Value = 0.3208215955797718
$ws.Range("Q6").Value = 7380.414011375183
$ws.Range("R6").Value = 44282.48406825109
$ws.Range("S6").Value = 0.04047947990550213
$ws.Range("T6").Value = 0.03722551003841511

$ws.Range("E7").Value = 2
$ws.Range("G7").Value = 164.7314145
$ws.Range("H7").Value = 329.462829
$ws.Range("I7").Value = 0.1559412251502966
$ws.Range("J7").Value = 0.1160318088037158
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 31.7643625
$ws.Range("N7").Value = 63.528725
$ws.Range("O7").Value = 0.1840389846252172
$ws.Range("P7").Value = 0.1516380468435678
$ws.Range("Q7").Value = 5232.588365315756
$ws.Range("R7").Value = 20930.35346126302
$ws.Range("S7").Value = 0.02869926473787297
$ws.Range("T7").Value = 0.01759483685872176

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 196.007169
$ws.Range("H8").Value = 588.021507
$ws.Range("I8").Value = 0.1855480824035615
$ws.Range("J8").Value = 0.2070922515896227
$ws.Range("K8").Value = 2
$ws.Range("M8").Value = 67.073376
$ws.Range("N8").Value = 134.146752
$ws.Range("O8").Value = 0.3886152607163267
$ws.Range("P8").Value = 0.3201976973989085
$ws.Range("Q8").Value = 13146.86254503254
$ws.Range("R8").Value = 78881.17527019526
$ws.Range("S8").Value = 0.07210681641867453
$ws.Range("T8").Value = 0.06631046210815263

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 196.007169
$ws.Range("H9").Value = 588.021507
$ws.Range("I9").Value = 0.1855480824035615
$ws.Range("J9").Value = 0.2070922515896227
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.219010666666667
$ws.Range("N9").Value = 12.657032
$ws.Range("O9").Value = 0.02444445215030819
$ws.Range("P9").Value = 0.03021133528677833
$ws.Range("Q9").Value = 826.956336754136
$ws.Range("R9").Value = 7442.607030787225
$ws.Range("S9").Value = 0.004535621221895301
$ws.Range("T9").Value = 0.006256533448067944

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 196.007169
$ws.Range("H10").Value = 588.021507
$ws.Range("I10").Value = 0.1855480824035615
$ws.Range("J10").Value = 0.2070922515896227
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 12.54988666666667
$ws.Range("N10").Value = 37.64966
$ws.Range("O10").Value = 0.07271256897710082
$ws.Range("P10").Value = 0.0898667635266472
$ws.Range("Q10").Value = 2459.86775680418
$ws.Range("R10").Value = 22138.80981123762
$ws.Range("S10").Value = 0.01349167774033775
$ws.Range("T10").Value = 0.01861071040180555

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 196.007169
$ws.Range("H11").Value = 588.021507
$ws.Range("I11").Value = 0.1855480824035615
$ws.Range("J11").Value = 0.2070922515896227
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 12.18648933333333
$ws.Range("N11").Value = 36.559468
$ws.Range("O11").Value = 0.07060708751994335
$ws.Range("P11").Value = 0.0872645613643264
$ws.Range("Q11").Value = 2388.639274275364
$ws.Range("R11").Value = 21497.75346847827
$ws.Range("S11").Value = 0.01310100969342593
$ws.Range("T11").Value = 0.01807181449691915

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 196.007169
$ws.Range("H12").Value = 588.021507
$ws.Range("I12").Value = 0.1855480824035615
$ws.Range("J12").Value = 0.2070922515896227
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 44.80271133333334
$ws.Range("N12").Value = 134.408134
$ws.Range("O12").Value = 0.2595816460111038
$ws.Range("P12").Value = 0.3208215955797718
$ws.Range("Q12").Value = 8781.652611970883
$ws.Range("R12").Value = 79034.87350773795
$ws.Range("S12").Value = 0.04816487664452041
$ws.Range("T12").Value = 0.06643966658719029

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 196.007169
$ws.Range("H13").Value = 588.021507
$ws.Range("I13").Value = 0.1855480824035615
$ws.Range("J13").Value = 0.2070922515896227
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 31.7643625
$ws.Range("N13").Value = 63.528725
$ws.Range("O13").Value = 0.1840389846252172
$ws.Range("P13").Value = 0.1516380468435678
$ws.Range("Q13").Value = 6226.042768714763
$ws.Range("R13").Value = 37356.25661228858
$ws.Range("S13").Value = 0.03414808068470759
$ws.Range("T13").Value = 0.03140306454748713

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 185.9706316666667
$ws.Range("H14").Value = 557.911895
$ws.Range("I14").Value = 0.1760471020788482
$ws.Range("J14").Value = 0.1964881031539942
$ws.Range("K14").Value = 2
$ws.Range("M14").Value = 67.073376
$ws.Range("N14").Value = 134.146752
$ws.Range("O14").Value = 0.3886152607163267
$ws.Range("P14").Value = 0.3201976973989085
$ws.Range("Q14").Value = 12473.67810273584
$ws.Range("R14").Value = 74842.06861641502
$ws.Range("S14").Value = 0.06841459047272537
$ws.Range("T14").Value = 0.06291503819618816

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 185.9706316666667
$ws.Range("H15").Value = 557.911895
$ws.Range("I15").Value = 0.1760471020788482
$ws.Range("J15").Value = 0.1964881031539942
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 4.219010666666667
$ws.Range("N15").Value = 12.657032
$ws.Range("O15").Value = 0.02444445215030819
$ws.Range("P15").Value = 0.03021133528677833
$ws.Range("Q15").Value = 784.6120786884044
$ws.Range("R15").Value = 7061.50870819564
$ws.Range("S15").Value = 0.004303374962966827
$ws.Range("T15").Value = 0.005936167964248406

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 185.9706316666667
$ws.Range("H16").Value = 557.911895
$ws.Range("I16").Value = 0.1760471020788482
$ws.Range("J16").Value = 0.1964881031539942
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 12.54988666666667
$ws.Range("N16").Value = 37.64966
$ws.Range("O16").Value = 0.07271256897710082
$ws.Range("P16").Value = 0.0898667635266472
$ws.Range("Q16").Value = 2333.910350745078
$ws.Range("R16").Value = 21005.1931567057
$ws.Range("S16").Value = 0.01280083705312696
$ws.Range("T16").Value = 0.01765774990193946

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 185.9706316666667
$ws.Range("H17").Value = 557.911895
$ws.Range("I17").Value = 0.1760471020788482
$ws.Range("J17").Value = 0.1964881031539942
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 12.18648933333333
$ws.Range("N17").Value = 36.559468
$ws.Range("O17").Value = 0.07060708751994335
$ws.Range("P17").Value = 0.0872645613643264
$ws.Range("Q17").Value = 2266.329119119095
$ws.Range("R17").Value = 20396.96207207186
$ws.Range("S17").Value = 0.01243017314411364
$ws.Range("T17").Value = 0.01714644813504183

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 185.9706316666667
$ws.Range("H18").Value = 557.911895
$ws.Range("I18").Value = 0.1760471020788482
$ws.Range("J18").Value = 0.1964881031539942
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 44.80271133333334
$ws.Range("N18").Value = 134.408134
$ws.Range("O18").Value = 0.2595816460111038
$ws.Range("P18").Value = 0.3208215955797718
$ws.Range("Q18").Value = 8331.988527039326
$ws.Range("R18").Value = 74987.89674335394
$ws.Range("S18").Value = 0.04569859653311222
$ws.Range("T18").Value = 0.06303762676630723

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 185.9706316666667
$ws.Range("H19").Value = 557.911895
$ws.Range("I19").Value = 0.1760471020788482
$ws.Range("J19").Value = 0.1964881031539942
$ws.Range("K19").Value = 2
$ws.Range("M19").Value = 31.7643625
$ws.Range("N19").Value = 63.528725
$ws.Range("O19").Value = 0.1840389846252172
$ws.Range("P19").Value = 0.1516380468435678
$ws.Range("Q19").Value = 5907.238558613979
$ws.Range("R19").Value = 35443.43135168387
$ws.Range("S19").Value = 0.03239952991280318
$ws.Range("T19").Value = 0.02979507219026916

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 171.1876676666667
$ws.Range("H20").Value = 513.563003
$ws.Range("I20").Value = 0.1620529679028636
$ws.Range("J20").Value = 0.1808690963822147
$ws.Range("K20").Value = 2
$ws.Range("M20").Value = 67.073376
$ws.Range("N20").Value = 134.146752
$ws.Range("O20").Value = 0.3886152607163267
$ws.Range("P20").Value = 0.3201976973989085
$ws.Range("Q20").Value = 11482.13479996938
$ws.Range("R20").Value = 68892.80879981625
$ws.Range("S20").Value = 0.06297625637142587
$ws.Range("T20").Value = 0.0579138681922064

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 171.1876676666667
$ws.Range("H21").Value = 513.563003
$ws.Range("I21").Value = 0.1620529679028636
$ws.Range("J21").Value = 0.1808690963822147
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 4.219010666666667
$ws.Range("N21").Value = 12.657032
$ws.Range("O21").Value = 0.02444445215030819
$ws.Range("P21").Value = 0.03021133528677833
$ws.Range("Q21").Value = 722.2425958874551
$ws.Range("R21").Value = 6500.183362987096
$ws.Range("S21").Value = 0.003961296019716979
$ws.Range("T21").Value = 0.005464296913819714

$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 171.1876676666667
$ws.Range("H22").Value = 513.563003
$ws.Range("I22").Value = 0.1620529679028636
$ws.Range("J22").Value = 0.1808690963822147
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 12.54988666666667
$ws.Range("N22").Value = 37.64966
$ws.Range("O22").Value = 0.07271256897710082
$ws.Range("P22").Value = 0.0898667635266472
$ws.Range("Q22").Value = 2148.385827947664
$ws.Range("R22").Value = 19335.47245152898
$ws.Range("S22").Value = 0.01178328760658088
$ws.Range("T22").Value = 0.01625412031385885

$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 171.1876676666667
$ws.Range("H23").Value = 513.563003
$ws.Range("I23").Value = 0.1620529679028636
$ws.Range("J23").Value = 0.1808690963822147
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 12.18648933333333
$ws.Range("N23").Value = 36.559468
$ws.Range("O23").Value = 0.07060708751994335
$ws.Range("P23").Value = 0.0872645613643264
$ws.Range("Q23").Value = 2086.176686018045
$ws.Range("R23").Value = 18775.5901741624
$ws.Range("S23").Value = 0.01144208808758406
$ws.Range("T23").Value = 0.01578346236015604

$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 171.1876676666667
$ws.Range("H24").Value = 513.563003
$ws.Range("I24").Value = 0.1620529679028636
$ws.Range("J24").Value = 0.1808690963822147
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 44.80271133333334
$ws.Range("N24").Value = 134.408134
$ws.Range("O24").Value = 0.2595816460111038
$ws.Range("P24").Value = 0.3208215955797718
$ws.Range("Q24").Value = 7669.671658296268
$ws.Range("R24").Value = 69027.0449246664
$ws.Range("S24").Value = 0.04206597614920991
$ws.Range("T24").Value = 0.05802671209241365

$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 171.1876676666667
$ws.Range("H25").Value = 513.563003
$ws.Range("I25").Value = 0.1620529679028636
$ws.Range("J25").Value = 0.1808690963822147
$ws.Range("K25").Value = 2
$ws.Range("M25").Value = 31.7643625
$ws.Range("N25").Value = 63.528725
$ws.Range("O25").Value = 0.1840389846252172
$ws.Range("P25").Value = 0.1516380468435678
$ws.Range("Q25").Value = 5437.66713129353
$ws.Range("R25").Value = 32626.00278776118
$ws.Range("S25").Value = 0.02982406366834593
$ws.Range("T25").Value = 0.02742663650976005

$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 173.5155153333334
$ws.Range("H26").Value = 520.546546
$ws.Range("I26").Value = 0.1642565999071482
$ws.Range("J26").Value = 0.1833285942521505
$ws.Range("K26").Value = 2
$ws.Range("M26").Value = 67.073376
$ws.Range("N26").Value = 134.146752
$ws.Range("O26").Value = 0.3886152607163267
$ws.Range("P26").Value = 0.3201976973989085
$ws.Range("Q26").Value = 11638.27140178643
$ws.Range("R26").Value = 69829.6284107186
$ws.Range("S26").Value = 0.06383262139729376
$ws.Range("T26").Value = 0.05870139374691737

$ws.Range("E27").Value = 3
$ws.Range("G27").Value = 173.5155153333334
$ws.Range("H27").Value = 520.546546
$ws.Range("I27").Value = 0.1642565999071482
$ws.Range("J27").Value = 0.1833285942521505
$ws.Range("K27").Value = 3
$ws.Range("M27").Value = 4.219010666666667
$ws.Range("N27").Value = 12.657032
$ws.Range("O27").Value = 0.02444445215030819
$ws.Range("P27").Value = 0.03021133528677833
$ws.Range("Q27").Value = 732.063810023497
$ws.Range("R27").Value = 6588.574290211473
$ws.Range("S27").Value = 0.0040151625968026
$ws.Range("T27").Value = 0.005538601628605463

$ws.Range("E28").Value = 3
$ws.Range("G28").Value = 173.5155153333334
$ws.Range("H28").Value = 520.546546
$ws.Range("I28").Value = 0.1642565999071482
$ws.Range("J28").Value = 0.1833285942521505
$ws.Range("K28").Value = 3
$ws.Range("M28").Value = 12.54988666666667
$ws.Range("N28").Value = 37.64966
$ws.Range("O28").Value = 0.07271256897710082
$ws.Range("P28").Value = 0.0898667635266472
$ws.Range("Q28").Value = 2177.600052341596
$ws.Range("R28").Value = 19598.40047107436
$ws.Range("S28").Value = 0.01194351935069256
$ws.Range("T28").Value = 0.01647514742733067

$ws.Range("E29").Value = 3
$ws.Range("G29").Value = 173.5155153333334
$ws.Range("H29").Value = 520.546546
$ws.Range("I29").Value = 0.1642565999071482
$ws.Range("J29").Value = 0.1833285942521505
$ws.Range("K29").Value = 3
$ws.Range("M29").Value = 12.18648933333333
$ws.Range("N29").Value = 36.559468
$ws.Range("O29").Value = 0.07060708751994335
$ws.Range("P29").Value = 0.0872645613643264
$ws.Range("Q29").Value = 2114.544976777503
$ws.Range("R29").Value = 19030.90479099753
$ws.Range("S29").Value = 0.01159768012537233
$ws.Range("T29").Value = 0.01599808936295249

$ws.Range("E30").Value = 3
$ws.Range("G30").Value = 173.5155153333334
$ws.Range("H30").Value = 520.546546
$ws.Range("I30").Value = 0.1642565999071482
$ws.Range("J30").Value = 0.1833285942521505
$ws.Range("K30").Value = 3
$ws.Range("M30").Value = 44.80271133333334
$ws.Range("N30").Value = 134.408134
$ws.Range("O30").Value = 0.2595816460111038
$ws.Range("P30").Value = 0.3208215955797718
$ws.Range("Q30").Value = 7773.965545333909
$ws.Range("R30").Value = 69965.68990800518
$ws.Range("S30").Value = 0.04263799857208483
$ws.Range("T30").Value = 0.05881577212337152

$ws.Range("E31").Value = 3
$ws.Range("G31").Value = 173.5155153333334
$ws.Range("H31").Value = 520.546546
$ws.Range("I31").Value = 0.1642565999071482
$ws.Range("J31").Value = 0.1833285942521505
$ws.Range("K31").Value = 2
$ws.Range("M31").Value = 31.7643625
$ws.Range("N31").Value = 63.528725
$ws.Range("O31").Value = 0.1840389846252172
$ws.Range("P31").Value = 0.1516380468435678
$ws.Range("Q31").Value = 5511.609728422309
$ws.Range("R31").Value = 33069.65837053386
$ws.Range("S31").Value = 0.03022961786490209
$ws.Range("T31").Value = 0.02779958996297304

$ws.Range("E32").Value = 2
$ws.Range("G32").Value = 164.956207
$ws.Range("H32").Value = 329.912414
$ws.Range("I32").Value = 0.1561540225572818
$ws.Range("J32").Value = 0.1161901458183021
$ws.Range("K32").Value = 2
$ws.Range("M32").Value = 67.073376
$ws.Range("N32").Value = 134.146752
$ws.Range("O32").Value = 0.3886152607163267
$ws.Range("P32").Value = 0.3201976973989085
$ws.Range("Q32").Value = 11064.16969564483
$ws.Range("R32").Value = 44256.67878257933
$ws.Range("S32").Value = 0.06068383618800124
$ws.Range("T32").Value = 0.03720381715146375

$ws.Range("E33").Value = 2
$ws.Range("G33").Value = 164.956207
$ws.Range("H33").Value = 329.912414
$ws.Range("I33").Value = 0.1561540225572818
$ws.Range("J33").Value = 0.1161901458183021
$ws.Range("K33").Value = 3
$ws.Range("M33").Value = 4.219010666666667
$ws.Range("N33").Value = 12.657032
$ws.Range("O33").Value = 0.02444445215030819
$ws.Range("P33").Value = 0.03021133528677833
$ws.Range("Q33").Value = 695.9519968658747
$ws.Range("R33").Value = 4175.711981195249
$ws.Range("S33").Value = 0.003817099532479622
$ws.Range("T33").Value = 0.00351025945233639

$ws.Range("E34").Value = 2
$ws.Range("G34").Value = 164.956207
$ws.Range("H34").Value = 329.912414
$ws.Range("I34").Value = 0.1561540225572818
$ws.Range("J34").Value = 0.1161901458183021
$ws.Range("K34").Value = 3
$ws.Range("M34").Value = 12.54988666666667
$ws.Range("N34").Value = 37.64966
$ws.Range("O34").Value = 0.07271256897710082
$ws.Range("P34").Value = 0.0898667635266472
$ws.Range("Q34").Value = 2070.181702813207
$ws.Range("R34").Value = 12421.09021687924
$ws.Range("S34").Value = 0.01135436013624811
$ws.Range("T34").Value = 0.01044163235838001

$ws.Range("E35").Value = 2
$ws.Range("G35").Value = 164.956207
$ws.Range("H35").Value = 329.912414
$ws.Range("I35").Value = 0.1561540225572818
$ws.Range("J35").Value = 0.1161901458183021
$ws.Range("K35").Value = 3
$ws.Range("M35").Value = 12.18648933333333
$ws.Range("N35").Value = 36.559468
$ws.Range("O35").Value = 0.07060708751994335
$ws.Range("P35").Value = 0.0872645613643264
$ws.Range("Q35").Value = 2010.237057072625
$ws.Range("R35").Value = 12061.42234243575
$ws.Range("S35").Value = 0.0110255807372932
$ws.Range("T35").Value = 0.01013928210969126

$ws.Range("E36").Value = 2
$ws.Range("G36").Value = 164.956207
$ws.Range("H36").Value = 329.912414
$ws.Range("I36").Value = 0.1561540225572818
$ws.Range("J36").Value = 0.1161901458183021
$ws.Range("K36").Value = 3
$ws.Range("M36").Value = 44.80271133333334
$ws.Range("N36").Value = 134.408134
$ws.Range("O36").Value = 0.2595816460111038
$ws.Range("P36").Value = 0.3208215955797718
$ws.Range("Q36").Value = 7390.485324862581
$ws.Range("R36").Value = 44342.91194917548
$ws.Range("S36").Value = 0.04053471820667424
$ws.Range("T36").Value = 0.03727630797207403

$ws.Range("E37").Value = 2
$ws.Range("G37").Value = 164.956207
$ws.Range("H37").Value = 329.912414
$ws.Range("I37").Value = 0.1561540225572818
$ws.Range("J37").Value = 0.1161901458183021
$ws.Range("K37").Value = 2
$ws.Range("M37").Value = 31.7643625
$ws.Range("N37").Value = 63.528725
$ws.Range("O37").Value = 0.1840389846252172
$ws.Range("P37").Value = 0.1516380468435678
$ws.Range("Q37").Value = 5239.728755773038
$ws.Range("R37").Value = 20958.91502309215
$ws.Range("S37").Value = 0.02873842775658541
$ws.Range("T37").Value = 0.01761884677435667
